$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.19972966816410173
$ws.Range("D2").Value = 0.19972966816410173
$ws.Range("E2").Value = 0.9794400430884299
$ws.Range("F2").Value = 0.0055290149473180596
$ws.Range("G2").Value = 0.4381

$ws.Range("C3").Value = 0.5437913689595433
$ws.Range("D3").Value = 0.5437913689595433
$ws.Range("E3").Value = 2.666659624183864
$ws.Range("F3").Value = 0.015053500237779193
$ws.Range("G3").Value = 0.0141

$ws.Range("C4").Value = 6.2113009675232265
$ws.Range("D4").Value = 2.0704336558410756
$ws.Range("E4").Value = 10.153051610853252
$ws.Range("F4").Value = 0.1719442895359476
$ws.Range("G4").Value = 0.0001

$ws.Range("C5").Value = 0.10629474801754313
$ws.Range("D5").Value = 0.10629474801754313
$ws.Range("E5").Value = 0.5212512168840034
$ws.Range("F5").Value = 0.0029425035149386655
$ws.Range("G5").Value = 0.8343

$ws.Range("C6").Value = 0.3631016987679727
$ws.Range("D6").Value = 0.12103389958932424
$ws.Range("E6").Value = 0.593529488726378
$ws.Range("F6").Value = 0.010051559882607036
$ws.Range("G6").Value = 0.9276

$ws.Range("C7").Value = 1.430894045606487
$ws.Range("D7").Value = 0.47696468186882895
$ws.Range("E7").Value = 2.3389530101128435
$ws.Range("F7").Value = 0.03961071301478601
$ws.Range("G7").Value = 0.0009

$ws.Range("C8").Value = 0.35105904465797444
$ws.Range("D8").Value = 0.11701968155265814
$ws.Range("E8").Value = 0.5738444518315675
$ws.Range("F8").Value = 0.009718189206174253
$ws.Range("G8").Value = 0.9372

$ws.Range("C9").Value = 26.917743851402953
$ws.Range("D9").Value = 0.20392230190456784
$ws.Range("F9").Value = 0.7451502296604492

$ws.Range("C10").Value = 36.1239153930998
